$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update resultado (G) and profit (H) for rows that were pending evaluation.
$updates = @(
    @{ Row = 23; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 24; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 28; Resultado = "Acierto"; Profit = 1.1 },
    @{ Row = 34; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 56; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 60; Resultado = "Acierto"; Profit = 0.91 },
    @{ Row = 65; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 66; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 71; Resultado = "Acierto"; Profit = 1.5 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 7).Value = $u.Resultado
    $ws.Cells.Item($u.Row, 8).Value = $u.Profit
}

# event_id on row 73 was stored as text; convert it to a numeric value.
$ws.Cells.Item(73, 1).Value = 14339215
